# Generate Report for Handoff
# Adds a new handoff row (bd1e060a-d5db-48a9-95cf-5645fc0d341e) to each of the
# three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-25 02:46:19"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8dd6ceb9f8441095549479ed65816d80d7203d78/e2e/bd1e060a-d5db-48a9-95cf-5645fc0d341e.md",
    "",
    "",
    "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-25 02:46:15"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J3").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8dd6ceb9f8441095549479ed65816d80d7203d78/e2e/bd1e060a-d5db-48a9-95cf-5645fc0d341e.md",
    "",
    "",
    "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f789b2d08e6505e0b3d7a054d52100f27d16a0c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.zh-cn.xlf",
    "",
    "",
    "bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-25 02:46:19"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("J3").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8dd6ceb9f8441095549479ed65816d80d7203d78/e2e/bd1e060a-d5db-48a9-95cf-5645fc0d341e.md",
    "",
    "",
    "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a7e6ad1e1d047ca0f9fe46fdece8e72fdae0c1c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.de-de.xlf",
    "",
    "",
    "bd1e060a-d5db-48a9-95cf-5645fc0d341e.7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5.de-de.xlf"
)
